$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 15:16"

# Rebuild the country data table (rows 4-202) to match the refreshed/re-sorted dataset
$data = New-Object 'object[,]' 199,8
$data[0,0] = "China"
$data[0,1] = 81218
$data[0,2] = 47
$data[0,3] = 73650
$data[0,4] = 4287
$data[0,5] = 1399
$data[0,6] = 4
$data[0,7] = 3281
$data[1,0] = "Italia"
$data[1,1] = 69176
$data[1,2] = 0
$data[1,3] = 8326
$data[1,4] = 54030
$data[1,5] = 3393
$data[1,6] = 0
$data[1,7] = 6820
$data[2,0] = "Estados Unidos"
$data[2,1] = 55081
$data[2,2] = 225
$data[2,3] = 379
$data[2,4] = 53917
$data[2,5] = 1175
$data[2,6] = 5
$data[2,7] = 785
$data[3,0] = "España"
$data[3,1] = 47610
$data[3,2] = 5552
$data[3,3] = 5367
$data[3,4] = 38809
$data[3,5] = 2636
$data[3,6] = 443
$data[3,7] = 3434
$data[4,0] = "Alemania"
$data[4,1] = 35531
$data[4,2] = 2540
$data[4,3] = 3540
$data[4,4] = 31810
$data[4,5] = 23
$data[4,6] = 22
$data[4,7] = 181
$data[5,0] = "Iran"
$data[5,1] = 27017
$data[5,2] = 2206
$data[5,3] = 9625
$data[5,4] = 15315
$data[5,5] = 0
$data[5,6] = 143
$data[5,7] = 2077
$data[6,0] = "Francia"
$data[6,1] = 22304
$data[6,2] = 0
$data[6,3] = 3281
$data[6,4] = 17923
$data[6,5] = 2516
$data[6,6] = 0
$data[6,7] = 1100
$data[7,0] = "Suiza"
$data[7,1] = 10537
$data[7,2] = 660
$data[7,3] = 131
$data[7,4] = 10257
$data[7,5] = 141
$data[7,6] = 27
$data[7,7] = 149
$data[8,0] = "Corea del Sur"
$data[8,1] = 9137
$data[8,2] = 100
$data[8,3] = 3730
$data[8,4] = 5281
$data[8,5] = 59
$data[8,6] = 6
$data[8,7] = 126
$data[9,0] = "Reino Unido"
$data[9,1] = 8227
$data[9,2] = 150
$data[9,3] = 135
$data[9,4] = 7659
$data[9,5] = 20
$data[9,6] = 11
$data[9,7] = 433
$data[10,0] = "Paises Bajos"
$data[10,1] = 6412
$data[10,2] = 852
$data[10,3] = 3
$data[10,4] = 6053
$data[10,5] = 582
$data[10,6] = 80
$data[10,7] = 356
$data[11,0] = "Austria"
$data[11,1] = 5588
$data[11,2] = 305
$data[11,3] = 9
$data[11,4] = 5549
$data[11,5] = 28
$data[11,6] = 2
$data[11,7] = 30
$data[12,0] = "Belgica"
$data[12,1] = 4937
$data[12,2] = 668
$data[12,3] = 547
$data[12,4] = 4212
$data[12,5] = 474
$data[12,6] = 56
$data[12,7] = 178
$data[13,0] = "Portugal"
$data[13,1] = 2995
$data[13,2] = 633
$data[13,3] = 22
$data[13,4] = 2930
$data[13,5] = 61
$data[13,6] = 10
$data[13,7] = 43
$data[14,0] = "Noruega"
$data[14,1] = 2971
$data[14,2] = 105
$data[14,3] = 6
$data[14,4] = 2951
$data[14,5] = 57
$data[14,6] = 2
$data[14,7] = 14
$data[15,0] = "Canada"
$data[15,1] = 2792
$data[15,2] = 0
$data[15,3] = 112
$data[15,4] = 2654
$data[15,5] = 1
$data[15,6] = 0
$data[15,7] = 26
$data[16,0] = "Suecia"
$data[16,1] = 2526
$data[16,2] = 227
$data[16,3] = 16
$data[16,4] = 2468
$data[16,5] = 158
$data[16,6] = 2
$data[16,7] = 42
$data[17,0] = "Australia"
$data[17,1] = 2431
$data[17,2] = 114
$data[17,3] = 118
$data[17,4] = 2304
$data[17,5] = 11
$data[17,6] = 1
$data[17,7] = 9
$data[18,0] = "Brasil"
$data[18,1] = 2274
$data[18,2] = 27
$data[18,3] = 2
$data[18,4] = 2225
$data[18,5] = 18
$data[18,6] = 1
$data[18,7] = 47
$data[19,0] = "Israel"
$data[19,1] = 2170
$data[19,2] = 240
$data[19,3] = 58
$data[19,4] = 2107
$data[19,5] = 37
$data[19,6] = 2
$data[19,7] = 5
$data[20,0] = "Turquia"
$data[20,1] = 1874
$data[20,2] = 2
$data[20,3] = 0
$data[20,4] = 1830
$data[20,5] = 0
$data[20,6] = 0
$data[20,7] = 44
$data[21,0] = "Malasia"
$data[21,1] = 1796
$data[21,2] = 172
$data[21,3] = 199
$data[21,4] = 1578
$data[21,5] = 64
$data[21,6] = 3
$data[21,7] = 19
$data[22,0] = "Dinamarca"
$data[22,1] = 1715
$data[22,2] = 124
$data[22,3] = 1
$data[22,4] = 1680
$data[22,5] = 69
$data[22,6] = 2
$data[22,7] = 34
$data[23,0] = "Chequia"
$data[23,1] = 1497
$data[23,2] = 103
$data[23,3] = 10
$data[23,4] = 1482
$data[23,5] = 19
$data[23,6] = 2
$data[23,7] = 5
$data[24,0] = "Irlanda"
$data[24,1] = 1329
$data[24,2] = 0
$data[24,3] = 5
$data[24,4] = 1317
$data[24,5] = 29
$data[24,6] = 0
$data[24,7] = 7
$data[25,0] = "Japon"
$data[25,1] = 1193
$data[25,2] = 0
$data[25,3] = 285
$data[25,4] = 865
$data[25,5] = 54
$data[25,6] = 0
$data[25,7] = 43
$data[26,0] = "Chile"
$data[26,1] = 1142
$data[26,2] = 220
$data[26,3] = 22
$data[26,4] = 1118
$data[26,5] = 7
$data[26,6] = 0
$data[26,7] = 2
$data[27,0] = "Luxemburgo"
$data[27,1] = 1099
$data[27,2] = 0
$data[27,3] = 6
$data[27,4] = 1085
$data[27,5] = 3
$data[27,6] = 0
$data[27,7] = 8
$data[28,0] = "Ecuador"
$data[28,1] = 1082
$data[28,2] = 0
$data[28,3] = 3
$data[28,4] = 1052
$data[28,5] = 2
$data[28,6] = 0
$data[28,7] = 27
$data[29,0] = "Pakistan"
$data[29,1] = 1022
$data[29,2] = 50
$data[29,3] = 21
$data[29,4] = 993
$data[29,5] = 5
$data[29,6] = 1
$data[29,7] = 8
$data[30,0] = "Polonia"
$data[30,1] = 957
$data[30,2] = 56
$data[30,3] = 2
$data[30,4] = 942
$data[30,5] = 3
$data[30,6] = 3
$data[30,7] = 13
$data[31,0] = "Tailandia"
$data[31,1] = 934
$data[31,2] = 107
$data[31,3] = 70
$data[31,4] = 860
$data[31,5] = 11
$data[31,6] = 0
$data[31,7] = 4
$data[32,0] = "Rumania"
$data[32,1] = 906
$data[32,2] = 112
$data[32,3] = 86
$data[32,4] = 807
$data[32,5] = 18
$data[32,6] = 1
$data[32,7] = 13
$data[33,0] = "Arabia Saudita"
$data[33,1] = 900
$data[33,2] = 133
$data[33,3] = 29
$data[33,4] = 869
$data[33,5] = 0
$data[33,6] = 1
$data[33,7] = 2
$data[34,0] = "Finlandia"
$data[34,1] = 880
$data[34,2] = 88
$data[34,3] = 10
$data[34,4] = 867
$data[34,5] = 22
$data[34,6] = 2
$data[34,7] = 3
$data[35,0] = "Indonesia"
$data[35,1] = 790
$data[35,2] = 104
$data[35,3] = 31
$data[35,4] = 701
$data[35,5] = 0
$data[35,6] = 3
$data[35,7] = 58
$data[36,0] = "Grecia"
$data[36,1] = 743
$data[36,2] = 0
$data[36,3] = 29
$data[36,4] = 694
$data[36,5] = 35
$data[36,6] = 0
$data[36,7] = 20
$data[37,0] = "Islandia"
$data[37,1] = 737
$data[37,2] = 89
$data[37,3] = 56
$data[37,4] = 679
$data[37,5] = 11
$data[37,6] = 0
$data[37,7] = 2
$data[38,0] = "Crucero"
$data[38,1] = 712
$data[38,2] = 0
$data[38,3] = 587
$data[38,4] = 115
$data[38,5] = 15
$data[38,6] = 0
$data[38,7] = 10
$data[39,0] = "Sudafrica"
$data[39,1] = 709
$data[39,2] = 155
$data[39,3] = 12
$data[39,4] = 697
$data[39,5] = 2
$data[39,6] = 0
$data[39,7] = 0
$data[40,0] = "Rusia"
$data[40,1] = 658
$data[40,2] = 163
$data[40,3] = 29
$data[40,4] = 628
$data[40,5] = 8
$data[40,6] = 0
$data[40,7] = 1
$data[41,0] = "Filipinas"
$data[41,1] = 636
$data[41,2] = 84
$data[41,3] = 26
$data[41,4] = 572
$data[41,5] = 1
$data[41,6] = 3
$data[41,7] = 38
$data[42,0] = "India"
$data[42,1] = 606
$data[42,2] = 70
$data[42,3] = 42
$data[42,4] = 554
$data[42,5] = 0
$data[42,6] = 0
$data[42,7] = 10
$data[43,0] = "Singapur"
$data[43,1] = 558
$data[43,2] = 0
$data[43,3] = 156
$data[43,4] = 400
$data[43,5] = 17
$data[43,6] = 0
$data[43,7] = 2
$data[44,0] = "Eslovenia"
$data[44,1] = 528
$data[44,2] = 48
$data[44,3] = 10
$data[44,4] = 513
$data[44,5] = 14
$data[44,6] = 1
$data[44,7] = 5
$data[45,0] = "Catar"
$data[45,1] = 526
$data[45,2] = 0
$data[45,3] = 41
$data[45,4] = 485
$data[45,5] = 6
$data[45,6] = 0
$data[45,7] = 0
$data[46,0] = "Panama"
$data[46,1] = 443
$data[46,2] = 0
$data[46,3] = 1
$data[46,4] = 436
$data[46,5] = 33
$data[46,6] = 0
$data[46,7] = 6
$data[47,0] = "Egipto"
$data[47,1] = 442
$data[47,2] = 40
$data[47,3] = 80
$data[47,4] = 341
$data[47,5] = 0
$data[47,6] = 1
$data[47,7] = 21
$data[48,0] = "Barein"
$data[48,1] = 419
$data[48,2] = 27
$data[48,3] = 177
$data[48,4] = 239
$data[48,5] = 2
$data[48,6] = 0
$data[48,7] = 3
$data[49,0] = "Croacia"
$data[49,1] = 418
$data[49,2] = 36
$data[49,3] = 16
$data[49,4] = 401
$data[49,5] = 6
$data[49,6] = 0
$data[49,7] = 1
$data[50,0] = "Peru"
$data[50,1] = 416
$data[50,2] = 0
$data[50,3] = 1
$data[50,4] = 408
$data[50,5] = 9
$data[50,6] = 0
$data[50,7] = 7
$data[51,0] = "Hong Kong"
$data[51,1] = 410
$data[51,2] = 23
$data[51,3] = 102
$data[51,4] = 304
$data[51,5] = 4
$data[51,6] = 0
$data[51,7] = 4
$data[52,0] = "Mexico"
$data[52,1] = 405
$data[52,2] = 38
$data[52,3] = 4
$data[52,4] = 396
$data[52,5] = 1
$data[52,6] = 1
$data[52,7] = 5
$data[53,0] = "Estonia"
$data[53,1] = 404
$data[53,2] = 35
$data[53,3] = 8
$data[53,4] = 396
$data[53,5] = 5
$data[53,6] = 0
$data[53,7] = 0
$data[54,0] = "Argentina"
$data[54,1] = 387
$data[54,2] = 0
$data[54,3] = 52
$data[54,4] = 328
$data[54,5] = 0
$data[54,6] = 1
$data[54,7] = 7
$data[55,0] = "Serbia"
$data[55,1] = 384
$data[55,2] = 81
$data[55,3] = 15
$data[55,4] = 365
$data[55,5] = 21
$data[55,6] = 1
$data[55,7] = 4
$data[56,0] = "Colombia"
$data[56,1] = 378
$data[56,2] = 0
$data[56,3] = 6
$data[56,4] = 369
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 3
$data[57,0] = "Irak"
$data[57,1] = 346
$data[57,2] = 30
$data[57,3] = 103
$data[57,4] = 214
$data[57,5] = 0
$data[57,6] = 2
$data[57,7] = 29
$data[58,0] = "Libano"
$data[58,1] = 333
$data[58,2] = 15
$data[58,3] = 8
$data[58,4] = 321
$data[58,5] = 4
$data[58,6] = 0
$data[58,7] = 4
$data[59,0] = "Emiratos Arabes Unidos"
$data[59,1] = 333
$data[59,2] = 85
$data[59,3] = 52
$data[59,4] = 279
$data[59,5] = 2
$data[59,6] = 0
$data[59,7] = 2
$data[60,0] = "Republica Dominicana"
$data[60,1] = 312
$data[60,2] = 0
$data[60,3] = 3
$data[60,4] = 303
$data[60,5] = 0
$data[60,6] = 0
$data[60,7] = 6
$data[61,0] = "Armenia"
$data[61,1] = 265
$data[61,2] = 16
$data[61,3] = 16
$data[61,4] = 249
$data[61,5] = 6
$data[61,6] = 0
$data[61,7] = 0
$data[62,0] = "Argelia"
$data[62,1] = 264
$data[62,2] = 0
$data[62,3] = 65
$data[62,4] = 180
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 19
$data[63,0] = "Lituania"
$data[63,1] = 255
$data[63,2] = 46
$data[63,3] = 1
$data[63,4] = 250
$data[63,5] = 1
$data[63,6] = 2
$data[63,7] = 4
$data[64,0] = "Taiwan"
$data[64,1] = 235
$data[64,2] = 19
$data[64,3] = 29
$data[64,4] = 204
$data[64,5] = 0
$data[64,6] = 0
$data[64,7] = 2
$data[65,0] = "Hungria"
$data[65,1] = 226
$data[65,2] = 39
$data[65,3] = 21
$data[65,4] = 195
$data[65,5] = 6
$data[65,6] = 1
$data[65,7] = 10
$data[66,0] = "Letonia"
$data[66,1] = 221
$data[66,2] = 24
$data[66,3] = 1
$data[66,4] = 220
$data[66,5] = 0
$data[66,6] = 0
$data[66,7] = 0
$data[67,0] = "Bulgaria"
$data[67,1] = 220
$data[67,2] = 2
$data[67,3] = 4
$data[67,4] = 213
$data[67,5] = 8
$data[67,6] = 0
$data[67,7] = 3
$data[68,0] = "Eslovaquia"
$data[68,1] = 216
$data[68,2] = 12
$data[68,3] = 7
$data[68,4] = 209
$data[68,5] = 2
$data[68,6] = 0
$data[68,7] = 0
$data[69,0] = "Nueva Zelanda"
$data[69,1] = 205
$data[69,2] = 0
$data[69,3] = 22
$data[69,4] = 183
$data[69,5] = 0
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = "Kuwait"
$data[70,1] = 195
$data[70,2] = 4
$data[70,3] = 43
$data[70,4] = 152
$data[70,5] = 6
$data[70,6] = 0
$data[70,7] = 0
$data[71,0] = "Uruguay"
$data[71,1] = 189
$data[71,2] = 0
$data[71,3] = 0
$data[71,4] = 189
$data[71,5] = 3
$data[71,6] = 0
$data[71,7] = 0
$data[72,0] = "Principado de Andorra"
$data[72,1] = 188
$data[72,2] = 24
$data[72,3] = 1
$data[72,4] = 186
$data[72,5] = 6
$data[72,6] = 0
$data[72,7] = 1
$data[73,0] = "San Marino"
$data[73,1] = 187
$data[73,2] = 0
$data[73,3] = 4
$data[73,4] = 162
$data[73,5] = 12
$data[73,6] = 0
$data[73,7] = 21
$data[74,0] = "Republica de Macedonia"
$data[74,1] = 177
$data[74,2] = 29
$data[74,3] = 1
$data[74,4] = 174
$data[74,5] = 1
$data[74,6] = 0
$data[74,7] = 2
$data[75,0] = "Costa Rica"
$data[75,1] = 177
$data[75,2] = 0
$data[75,3] = 2
$data[75,4] = 173
$data[75,5] = 4
$data[75,6] = 0
$data[75,7] = 2
$data[76,0] = "Marruecos"
$data[76,1] = 170
$data[76,2] = 0
$data[76,3] = 6
$data[76,4] = 159
$data[76,5] = 1
$data[76,6] = 0
$data[76,7] = 5
$data[77,0] = "Bosnia y Herzegovina"
$data[77,1] = 168
$data[77,2] = 0
$data[77,3] = 2
$data[77,4] = 163
$data[77,5] = 1
$data[77,6] = 0
$data[77,7] = 3
$data[78,0] = "Jordania"
$data[78,1] = 154
$data[78,2] = 0
$data[78,3] = 1
$data[78,4] = 153
$data[78,5] = 0
$data[78,6] = 0
$data[78,7] = 0
$data[79,0] = "Albania"
$data[79,1] = 146
$data[79,2] = 23
$data[79,3] = 17
$data[79,4] = 124
$data[79,5] = 3
$data[79,6] = 0
$data[79,7] = 5
$data[80,0] = "Vietnam"
$data[80,1] = 141
$data[80,2] = 7
$data[80,3] = 17
$data[80,4] = 124
$data[80,5] = 3
$data[80,6] = 0
$data[80,7] = 0
$data[81,0] = "Islas Feroe"
$data[81,1] = 132
$data[81,2] = 10
$data[81,3] = 38
$data[81,4] = 94
$data[81,5] = 2
$data[81,6] = 0
$data[81,7] = 0
$data[82,0] = "Malta"
$data[82,1] = 129
$data[82,2] = 19
$data[82,3] = 2
$data[82,4] = 127
$data[82,5] = 1
$data[82,6] = 0
$data[82,7] = 0
$data[83,0] = "Moldavia"
$data[83,1] = 125
$data[83,2] = 0
$data[83,3] = 2
$data[83,4] = 122
$data[83,5] = 20
$data[83,6] = 0
$data[83,7] = 1
$data[84,0] = "Republica de Chipre"
$data[84,1] = 124
$data[84,2] = 0
$data[84,3] = 3
$data[84,4] = 118
$data[84,5] = 3
$data[84,6] = 0
$data[84,7] = 3
$data[85,0] = "Tunez"
$data[85,1] = 119
$data[85,2] = 5
$data[85,3] = 2
$data[85,4] = 113
$data[85,5] = 11
$data[85,6] = 0
$data[85,7] = 4
$data[86,0] = "Burkina Faso"
$data[86,1] = 114
$data[86,2] = 0
$data[86,3] = 7
$data[86,4] = 103
$data[86,5] = 0
$data[86,6] = 0
$data[86,7] = 4
$data[87,0] = "Ucrania"
$data[87,1] = 113
$data[87,2] = 11
$data[87,3] = 1
$data[87,4] = 108
$data[87,5] = 0
$data[87,6] = 1
$data[87,7] = 4
$data[88,0] = "Brunei"
$data[88,1] = 109
$data[88,2] = 5
$data[88,3] = 2
$data[88,4] = 107
$data[88,5] = 1
$data[88,6] = 0
$data[88,7] = 0
$data[89,0] = "Sri Lanka"
$data[89,1] = 102
$data[89,2] = 0
$data[89,3] = 3
$data[89,4] = 99
$data[89,5] = 2
$data[89,6] = 0
$data[89,7] = 0
$data[90,0] = "Senegal"
$data[90,1] = 99
$data[90,2] = 13
$data[90,3] = 9
$data[90,4] = 90
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 0
$data[91,0] = "Oman"
$data[91,1] = 99
$data[91,2] = 15
$data[91,3] = 17
$data[91,4] = 82
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 0
$data[92,0] = "Reunion"
$data[92,1] = 94
$data[92,2] = 0
$data[92,3] = 1
$data[92,4] = 93
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 0
$data[93,0] = "Camboya"
$data[93,1] = 93
$data[93,2] = 2
$data[93,3] = 6
$data[93,4] = 87
$data[93,5] = 1
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = "Venezuela"
$data[94,1] = 91
$data[94,2] = 7
$data[94,3] = 15
$data[94,4] = 76
$data[94,5] = 2
$data[94,6] = 0
$data[94,7] = 0
$data[95,0] = "Azerbaiyan"
$data[95,1] = 87
$data[95,2] = 0
$data[95,3] = 10
$data[95,4] = 76
$data[95,5] = 6
$data[95,6] = 0
$data[95,7] = 1
$data[96,0] = "Bielorrusia"
$data[96,1] = 86
$data[96,2] = 5
$data[96,3] = 29
$data[96,4] = 57
$data[96,5] = 2
$data[96,6] = 0
$data[96,7] = 0
$data[97,0] = "Kazajistan"
$data[97,1] = 80
$data[97,2] = 8
$data[97,3] = 0
$data[97,4] = 80
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 0
$data[98,0] = "Afganistan"
$data[98,1] = 79
$data[98,2] = 5
$data[98,3] = 2
$data[98,4] = 75
$data[98,5] = 0
$data[98,6] = 1
$data[98,7] = 2
$data[99,0] = "Guadalupe"
$data[99,1] = 73
$data[99,2] = 0
$data[99,3] = 0
$data[99,4] = 72
$data[99,5] = 4
$data[99,6] = 0
$data[99,7] = 1
$data[100,0] = "Costa de Marfil"
$data[100,1] = 73
$data[100,2] = 0
$data[100,3] = 2
$data[100,4] = 71
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 0
$data[101,0] = "Georgia"
$data[101,1] = 73
$data[101,2] = 3
$data[101,3] = 10
$data[101,4] = 63
$data[101,5] = 1
$data[101,6] = 0
$data[101,7] = 0
$data[102,0] = "Camerun"
$data[102,1] = 70
$data[102,2] = 4
$data[102,3] = 2
$data[102,4] = 67
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 1
$data[103,0] = "Ghana"
$data[103,1] = 68
$data[103,2] = 15
$data[103,3] = 0
$data[103,4] = 65
$data[103,5] = 0
$data[103,6] = 1
$data[103,7] = 3
$data[104,0] = "Estado de Palestina"
$data[104,1] = 62
$data[104,2] = 2
$data[104,3] = 16
$data[104,4] = 46
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 0
$data[105,0] = "Trinidad yTobago"
$data[105,1] = 57
$data[105,2] = 0
$data[105,3] = 0
$data[105,4] = 57
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 0
$data[106,0] = "Martinica"
$data[106,1] = 57
$data[106,2] = 0
$data[106,3] = 0
$data[106,4] = 56
$data[106,5] = 7
$data[106,6] = 0
$data[106,7] = 1
$data[107,0] = "Uzbekistan"
$data[107,1] = 56
$data[107,2] = 6
$data[107,3] = 0
$data[107,4] = 56
$data[107,5] = 4
$data[107,6] = 0
$data[107,7] = 0
$data[108,0] = "Montenegro"
$data[108,1] = 52
$data[108,2] = 5
$data[108,3] = 0
$data[108,4] = 51
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 1
$data[109,0] = "Liechtenstein"
$data[109,1] = 51
$data[109,2] = 0
$data[109,3] = 0
$data[109,4] = 51
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 0
$data[110,0] = "Consejo Danes para los Refugiados"
$data[110,1] = 48
$data[110,2] = 3
$data[110,3] = 0
$data[110,4] = 46
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 2
$data[111,0] = "Mauricio"
$data[111,1] = 48
$data[111,2] = 6
$data[111,3] = 0
$data[111,4] = 46
$data[111,5] = 1
$data[111,6] = 0
$data[111,7] = 2
$data[112,0] = "Cuba"
$data[112,1] = 48
$data[112,2] = 0
$data[112,3] = 1
$data[112,4] = 46
$data[112,5] = 2
$data[112,6] = 0
$data[112,7] = 1
$data[113,0] = "Nigeria"
$data[113,1] = 46
$data[113,2] = 2
$data[113,3] = 2
$data[113,4] = 43
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 1
$data[114,0] = "Kirguistan"
$data[114,1] = 44
$data[114,2] = 2
$data[114,3] = 0
$data[114,4] = 44
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 0
$data[115,0] = "Ruanda"
$data[115,1] = 40
$data[115,2] = 0
$data[115,3] = 0
$data[115,4] = 40
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 0
$data[116,0] = "Puerto Rico"
$data[116,1] = 39
$data[116,2] = 0
$data[116,3] = 1
$data[116,4] = 36
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 2
$data[117,0] = "Banglades"
$data[117,1] = 39
$data[117,2] = 0
$data[117,3] = 7
$data[117,4] = 27
$data[117,5] = 0
$data[117,6] = 1
$data[117,7] = 5
$data[118,0] = "Paraguay"
$data[118,1] = 37
$data[118,2] = 10
$data[118,3] = 0
$data[118,4] = 34
$data[118,5] = 1
$data[118,6] = 1
$data[118,7] = 3
$data[119,0] = "Honduras"
$data[119,1] = 36
$data[119,2] = 6
$data[119,3] = 0
$data[119,4] = 36
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = "Mayotte"
$data[120,1] = 36
$data[120,2] = 0
$data[120,3] = 0
$data[120,4] = 36
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 0
$data[121,0] = "Bolivia"
$data[121,1] = 32
$data[121,2] = 3
$data[121,3] = 0
$data[121,4] = 32
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 0
$data[122,0] = "Guam"
$data[122,1] = 32
$data[122,2] = 0
$data[122,3] = 0
$data[122,4] = 31
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 1
$data[123,0] = "Macao"
$data[123,1] = 30
$data[123,2] = 1
$data[123,3] = 10
$data[123,4] = 20
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = "Kenia"
$data[124,1] = 25
$data[124,2] = 0
$data[124,3] = 0
$data[124,4] = 25
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = "Polinesia Francesa"
$data[125,1] = 25
$data[125,2] = 0
$data[125,3] = 0
$data[125,4] = 25
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = "Jamaica"
$data[126,1] = 25
$data[126,2] = 4
$data[126,3] = 2
$data[126,4] = 22
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 1
$data[127,0] = "Isla de Man"
$data[127,1] = 23
$data[127,2] = 0
$data[127,3] = 0
$data[127,4] = 23
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = "Monaco"
$data[128,1] = 23
$data[128,2] = 0
$data[128,3] = 1
$data[128,4] = 22
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = "Togo"
$data[129,1] = 23
$data[129,2] = 3
$data[129,3] = 1
$data[129,4] = 22
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = "Guayana Francesa"
$data[130,1] = 23
$data[130,2] = 0
$data[130,3] = 6
$data[130,4] = 17
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = "Guatemala"
$data[131,1] = 21
$data[131,2] = 0
$data[131,3] = 0
$data[131,4] = 20
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 1
$data[132,0] = "Madagascar"
$data[132,1] = 19
$data[132,2] = 2
$data[132,3] = 0
$data[132,4] = 19
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = "Barbados"
$data[133,1] = 18
$data[133,2] = 0
$data[133,3] = 0
$data[133,4] = 18
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = "Islas Virgenes de los Estados Unidos"
$data[134,1] = 17
$data[134,2] = 0
$data[134,3] = 0
$data[134,4] = 17
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = "Aruba"
$data[135,1] = 17
$data[135,2] = 0
$data[135,3] = 1
$data[135,4] = 16
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = "Gibraltar"
$data[136,1] = 15
$data[136,2] = 0
$data[136,3] = 5
$data[136,4] = 10
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = "Nueva Caledonia"
$data[137,1] = 14
$data[137,2] = 4
$data[137,3] = 0
$data[137,4] = 14
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 0
$data[138,0] = "Uganda"
$data[138,1] = 14
$data[138,2] = 5
$data[138,3] = 0
$data[138,4] = 14
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = "Maldivas"
$data[139,1] = 13
$data[139,2] = 0
$data[139,3] = 8
$data[139,4] = 5
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = "Etiopia"
$data[140,1] = 12
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 12
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = "Zambia"
$data[141,1] = 12
$data[141,2] = 9
$data[141,3] = 0
$data[141,4] = 12
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = "Tanzania"
$data[142,1] = 12
$data[142,2] = 0
$data[142,3] = 0
$data[142,4] = 12
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = "Republica de Yibuti"
$data[143,1] = 11
$data[143,2] = 8
$data[143,3] = 0
$data[143,4] = 11
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = "Mongolia"
$data[144,1] = 10
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 10
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = "Guinea Ecuatorial"
$data[145,1] = 9
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 9
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = "El Salvador"
$data[146,1] = 9
$data[146,2] = 4
$data[146,3] = 0
$data[146,4] = 9
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = "San Martin (Parte Francesa)"
$data[147,1] = 8
$data[147,2] = 0
$data[147,3] = 0
$data[147,4] = 8
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 0
$data[148,0] = "Dominica"
$data[148,1] = 7
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 7
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = "Surinam"
$data[149,1] = 7
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 7
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = "Seychelles"
$data[150,1] = 7
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 7
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 0
$data[151,0] = "Haiti"
$data[151,1] = 7
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 7
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 0
$data[152,0] = "Niger"
$data[152,1] = 7
$data[152,2] = 4
$data[152,3] = 0
$data[152,4] = 6
$data[152,5] = 0
$data[152,6] = 1
$data[152,7] = 1
$data[153,0] = "Namibia"
$data[153,1] = 7
$data[153,2] = 0
$data[153,3] = 2
$data[153,4] = 5
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = "Benin"
$data[154,1] = 6
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 6
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = "Bermudas"
$data[155,1] = 6
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 6
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = "Gabon"
$data[156,1] = 6
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 5
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 1
$data[157,0] = "Curazao"
$data[157,1] = 6
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 5
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 1
$data[158,0] = "Islas Caimanes"
$data[158,1] = 6
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 5
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 1
$data[159,0] = "Fiyi"
$data[159,1] = 5
$data[159,2] = 1
$data[159,3] = 0
$data[159,4] = 5
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = "Guyana"
$data[160,1] = 5
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 4
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 1
$data[161,0] = "Bahamas"
$data[161,1] = 5
$data[161,2] = 0
$data[161,3] = 1
$data[161,4] = 4
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = "Groenlandia"
$data[162,1] = 5
$data[162,2] = 0
$data[162,3] = 2
$data[162,4] = 3
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = "Suazilandia"
$data[163,1] = 4
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 4
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = "Santa Sede"
$data[164,1] = 4
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 4
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = "Guinea"
$data[165,1] = 4
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 4
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = "Congo"
$data[166,1] = 4
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 4
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = "Nepal"
$data[167,1] = 4
$data[167,2] = 2
$data[167,3] = 1
$data[167,4] = 3
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = "Cabo Verde"
$data[168,1] = 4
$data[168,2] = 1
$data[168,3] = 0
$data[168,4] = 3
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 1
$data[169,0] = "Liberia"
$data[169,1] = 3
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 3
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = "Santa Lucia"
$data[170,1] = 3
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 3
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = "Mozambique"
$data[171,1] = 3
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 3
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = "Antigua y Barbuda"
$data[172,1] = 3
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 3
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = "Republica de Africa Central"
$data[173,1] = 3
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 3
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = "Laos"
$data[174,1] = 3
$data[174,2] = 1
$data[174,3] = 0
$data[174,4] = 3
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = "San Bartolome"
$data[175,1] = 3
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 3
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = "Republica del Chad"
$data[176,1] = 3
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 3
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = "Birmania"
$data[177,1] = 3
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 3
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "Angola"
$data[178,1] = 3
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 3
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = "Gambia"
$data[179,1] = 3
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 2
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 1
$data[180,0] = "Zimbabue"
$data[180,1] = 3
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 2
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 1
$data[181,0] = "Sudan"
$data[181,1] = 3
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 2
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 1
$data[182,0] = "Mauritania"
$data[182,1] = 2
$data[182,2] = 0
$data[182,3] = 0
$data[182,4] = 2
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = "Nicaragua"
$data[183,1] = 2
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 2
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = "San Martin (Parte Holandesa)"
$data[184,1] = 2
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 2
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = "Butan"
$data[185,1] = 2
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 2
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = "Mali"
$data[186,1] = 2
$data[186,2] = 2
$data[186,3] = 0
$data[186,4] = 2
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = "Guinea-Bisau"
$data[187,1] = 2
$data[187,2] = 2
$data[187,3] = 0
$data[187,4] = 2
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = "Eritrea"
$data[188,1] = 1
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 1
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0
$data[189,0] = "Siria"
$data[189,1] = 1
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 1
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = "Papua Nueva Guinea"
$data[190,1] = 1
$data[190,2] = 0
$data[190,3] = 0
$data[190,4] = 1
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 0
$data[191,0] = "Islas Turcas y Caicos"
$data[191,1] = 1
$data[191,2] = 0
$data[191,3] = 0
$data[191,4] = 1
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = "Somalia"
$data[192,1] = 1
$data[192,2] = 0
$data[192,3] = 0
$data[192,4] = 1
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = "Timor Oriental"
$data[193,1] = 1
$data[193,2] = 0
$data[193,3] = 0
$data[193,4] = 1
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = "Belice"
$data[194,1] = 1
$data[194,2] = 0
$data[194,3] = 0
$data[194,4] = 1
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 0
$data[195,0] = "Libia"
$data[195,1] = 1
$data[195,2] = 0
$data[195,3] = 0
$data[195,4] = 1
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0
$data[196,0] = "Granada"
$data[196,1] = 1
$data[196,2] = 0
$data[196,3] = 0
$data[196,4] = 1
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = "Montserrat"
$data[197,1] = 1
$data[197,2] = 0
$data[197,3] = 0
$data[197,4] = 1
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = "San Vicente y las Granadinas"
$data[198,1] = 1
$data[198,2] = 0
$data[198,3] = 0
$data[198,4] = 1
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0

$ws.Range("A4:H202").Value = $data

